$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Dyson Daniels"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Atlanta Hawks"

$ws.Range("A3").Value = "Luguentz Dort"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Oklahoma City Thunder"

$ws.Range("A4").Value = "Jamal Murray"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Denver Nuggets"

$ws.Range("A5").Value = "Josh Hart"
$ws.Range("B5").Value = "SF,PF"
$ws.Range("C5").Value = "New York Knicks"

$ws.Range("A6").Value = "Kyle Filipowski"
$ws.Range("B6").Value = "PF,C"
$ws.Range("C6").Value = "Utah Jazz"

$ws.Range("A7").Value = "Tari Eason"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Houston Rockets"

$ws.Range("A8").Value = "Alexandre Sarr"
$ws.Range("B8").Value = "PF,C"
$ws.Range("C8").Value = "Washington Wizards"

$ws.Range("A9").Value = "Naz Reid"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "Minnesota Timberwolves"

$ws.Range("A10").Value = "Myles Turner"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Indiana Pacers"

$ws.Range("A11").Value = "Victor Wembanyama"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "San Antonio Spurs"

$ws.Range("A12").Value = "De'Andre Hunter"
$ws.Range("B12").Value = "SF,PF"
$ws.Range("C12").Value = "Atlanta Hawks"

$ws.Range("A13").Value = "Donovan Mitchell"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("C13").Value = "Cleveland Cavaliers"

$ws.Range("A14").Value = "Kevin Huerter"
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "Sacramento Kings"

$ws.Range("A15").Value = "Ochai Agbaji"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Toronto Raptors"

$ws.Range("A16").Value = "Domantas Sabonis"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Sacramento Kings"

